$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.981.77'
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '3.396.11'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.47%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '3.399.14'
$ws.Range("E8").Value = '  -1.77%  '
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.56'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.395'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = '3.972.63'
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("E15").Value = '  +0.64%  '
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("D17").Value = '3.392.62'
$ws.Range("E17").Value = '  -1.86%  '
$ws.Range("D18").Value = '61.042.64'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '390.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.564'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.12%  '
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000119'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.36%  '
$ws.Range("D27").Value = '3.532.37'
$ws.Range("E27").Value = '  -1.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.180'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.45'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.48%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.74%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.90'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.78%  '
$ws.Range("D37").Value = '3.420.50'
$ws.Range("E37").Value = '  -1.67%  '
$ws.Range("E38").Value = '  -2.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '166.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0784'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.10%  '
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.49'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.83%  '
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("D48").Value = '2.544.78'
$ws.Range("E48").Value = '  -2.24%  '
$ws.Range("E49").Value = '  -3.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.82%  '
